$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.980.52"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 4).Value = "1.867.34"
$ws.Cells.Item(3, 5).Value = "  -2.96%  "
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$ws.Cells.Item(5, 4).Value = "317.99"
$ws.Cells.Item(5, 5).Value = "  -2.35%  "
$ws.Cells.Item(6, 5).Value = "  +0.10%  "
$ws.Cells.Item(7, 4).Value = "0.5079"
$ws.Cells.Item(7, 5).Value = "  -1.65%  "
$ws.Cells.Item(8, 5).Value = "  -2.40%  "
$ws.Cells.Item(9, 4).Value = "0.08147"
$ws.Cells.Item(9, 5).Value = "  -3.84%  "
$ws.Cells.Item(10, 4).Value = "41.93"
$ws.Cells.Item(11, 4).Value = "1.086"
$ws.Cells.Item(11, 5).Value = "  -3.24%  "
$ws.Cells.Item(12, 4).Value = "22.57"
$ws.Cells.Item(12, 5).Value = "  +6.45%  "
$ws.Cells.Item(13, 4).Value = "1.867.99"
$ws.Cells.Item(13, 5).Value = "  -2.73%  "
$ws.Cells.Item(15, 4).Value = "7.138"
$ws.Cells.Item(15, 5).Value = "  -2.82%  "
$ws.Cells.Item(16, 5).Value = "  +0.13%  "
$ws.Cells.Item(17, 4).Value = "91.56"
$ws.Cells.Item(17, 5).Value = "  -2.95%  "
$ws.Cells.Item(18, 4).Value = "0.00001073"
$ws.Cells.Item(18, 5).Value = "  -3.78%  "
$ws.Cells.Item(19, 4).Value = "0.06329"
$ws.Cells.Item(19, 5).Value = "  -6.59%  "
$ws.Cells.Item(20, 4).Value = "17.80"
$ws.Cells.Item(20, 5).Value = "  -1.19%  "
$ws.Cells.Item(21, 5).Value = "  +0.03%  "
$ws.Cells.Item(22, 4).Value = "29.960.14"
$ws.Cells.Item(22, 5).Value = "  -0.40%  "
$ws.Cells.Item(23, 4).Value = "5.780"
$ws.Cells.Item(23, 5).Value = "  -4.63%  "
$ws.Cells.Item(24, 5).Value = "  -1.31%  "
$ws.Cells.Item(25, 4).Value = "2.203"
$ws.Cells.Item(25, 5).Value = "  +0.12%  "
$ws.Cells.Item(26, 4).Value = "2.087.18"
$ws.Cells.Item(26, 5).Value = "  -2.53%  "
$ws.Cells.Item(27, 4).Value = "160.22"
$ws.Cells.Item(27, 5).Value = "  +0.06%  "
$ws.Cells.Item(28, 4).Value = "20.81"
$ws.Cells.Item(28, 5).Value = "  -0.84%  "
$ws.Cells.Item(29, 4).Value = "2.215"
$ws.Cells.Item(29, 5).Value = "  -10.22%  "
$ws.Cells.Item(30, 4).Value = "126.07"
$ws.Cells.Item(30, 5).Value = "  -2.30%  "
$ws.Cells.Item(31, 4).Value = "0.1031"
$ws.Cells.Item(31, 5).Value = "  -2.71%  "
$ws.Cells.Item(32, 4).Value = "1.037"
$ws.Cells.Item(32, 5).Value = "  -3.79%  "
$ws.Cells.Item(33, 4).Value = "5.846"
$ws.Cells.Item(33, 5).Value = "  -3.80%  "
$ws.Cells.Item(34, 4).Value = "3.736"
$ws.Cells.Item(34, 5).Value = "  +2.41%  "
$ws.Cells.Item(35, 5).Value = "  -3.65%  "
$ws.Cells.Item(36, 4).Value = "0.06319"
$ws.Cells.Item(36, 5).Value = "  -4.36%  "
$ws.Cells.Item(37, 4).Value = "5.156"
$ws.Cells.Item(37, 5).Value = "  -0.87%  "
$ws.Cells.Item(38, 4).Value = "0.2135"
$ws.Cells.Item(38, 5).Value = "  -4.09%  "
$ws.Cells.Item(39, 4).Value = "1.167"
$ws.Cells.Item(39, 5).Value = "  -6.27%  "
$ws.Cells.Item(40, 4).Value = "8.442"
$ws.Cells.Item(40, 5).Value = "  -6.33%  "
$ws.Cells.Item(41, 4).Value = "0.6245"
$ws.Cells.Item(41, 5).Value = "  -4.57%  "
$ws.Cells.Item(42, 4).Value = "1.206"
$ws.Cells.Item(42, 5).Value = "  -3.03%  "
$ws.Cells.Item(43, 5).Value = "  -1.78%  "
$ws.Cells.Item(44, 4).Value = "0.9999"
$ws.Cells.Item(44, 5).Value = "  -0.04%  "
$ws.Cells.Item(45, 4).Value = "0.5849"
$ws.Cells.Item(45, 5).Value = "  -4.73%  "
$ws.Cells.Item(46, 4).Value = "12.74"
$ws.Cells.Item(46, 5).Value = "  -3.00%  "
$ws.Cells.Item(47, 4).Value = "3.620"
$ws.Cells.Item(47, 5).Value = "  -3.49%  "
$ws.Cells.Item(49, 4).Value = "121.68"
$ws.Cells.Item(49, 5).Value = "  -3.15%  "
$ws.Cells.Item(50, 5).Value = "  -3.72%  "
$ws.Cells.Item(51, 4).Value = "1.151"
$ws.Cells.Item(51, 5).Value = "  +0.36%  "
